$d = $word.ActiveDocument

# --- 1. Merge "The system shall offer a way to set a custom time for the
#        clock," + " for hours, minutes and seconds" into a single run ---
$d.Content.Find.Execute(
    "The system shall offer a way to set a custom time for the clock, for hours, minutes and seconds",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The system shall offer a way to set a custom time for the clock, for hours, minutes and seconds",
    2) | Out-Null

# --- 2. Merge "T" + "he system shall offer a way to set a custom time
#        using the format: xx:xx:xx" into a single run ---
$d.Content.Find.Execute(
    "The system shall offer a way to set a custom time using the format: xx:xx:xx",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The system shall offer a way to set a custom time using the format: xx:xx:xx",
    2) | Out-Null

# --- 3. Update the "custom date" bullet (split into two runs) and insert
#        the two new bullets that follow it ---
$rng = $d.Content
$rng.Find.Execute("The system shall offer a way to set a custom date", $false) | Out-Null
# Extend the found range to include the paragraph mark so InsertXML
# replaces the whole paragraph (preserving its pPr exactly, including the
# empty <w:rPr/>) instead of merely appending a run inside it.
$pRng = $d.Range($rng.Start, $rng.End + 1)

$xml1 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Normal1"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:rPr></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="none"/>
    </w:rPr>
    <w:t xml:space="preserve">The system shall offer a way to set a custom date </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="none"/>
    </w:rPr>
    <w:t>for day, month and year</w:t>
  </w:r>
</w:p>
"@
$pRng.InsertXML($xml1)

# Re-locate the (now split) paragraph and collapse to its end so the two
# brand-new bullets can be inserted right after it, each as its own
# <w:p> with its own pPr/rPr, rather than being merged into the existing
# paragraph.
$rng2 = $d.Content
$rng2.Find.Execute("for day, month and year", $false) | Out-Null
$afterRng = $d.Range($rng2.End, $rng2.End)
$afterRng.Collapse(0)

$xml2 = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Normal1"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="none"/>
    </w:rPr>
    <w:t>The system shall use the current date for those segments (day, month and year) that weren&#8217;t set by a custom time flag</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Normal1"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="none"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="none"/>
    </w:rPr>
    <w:t>The system shall have the feature of hide the seconds of the clock</w:t>
  </w:r>
</w:p>
"@
$afterRng.InsertXML($xml2)

# --- 4. Merge "The system shall implement each clock digit " + "and colon"
#        + " on a separate Ncurses window" into a single run ---
$d.Content.Find.Execute(
    "The system shall implement each clock digit and colon on a separate Ncurses window",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The system shall implement each clock digit and colon on a separate Ncurses window",
    2) | Out-Null
